$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new log entry row (row 10), matching the date formatting used
# by the other rows in column A (reuse the existing date style, then set
# the serial date value so no time-of-day fraction gets introduced)
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A10").Value = 43453

$ws.Range("B10").Value = "Read some of the potential ""background""-section-related papers, wrote summaries for them. May need more material to read on."

$ws.Range("C10").Value = 1

# Update the selection to match the author's final cursor position
$ws.Range("C12").Select()

$wb.Save()
